# Lecture partielle de l'EDT M1 MIAGE.
# Shift every schedule date in column A forward by 3 years (+1096 days for
# this particular set of dates), and update the matching weekday label in
# column B to reflect the new day-of-week (vendredi->mardi, lundi->vendredi).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new Excel serial date (old date shifted by +3 years / +1096 days)
$dateRows = @{
    2  = 46028
    4  = 46031
    7  = 46049
    9  = 46052
    12 = 46056
    15 = 46091
    18 = 46098
    21 = 46105
    24 = 46112
    27 = 46147
    30 = 46154
    33 = 46161
}

foreach ($row in $dateRows.Keys) {
    $ws.Cells.Item($row, 1).Value = $dateRows[$row]
}

# row -> new weekday label for column B
$dayRows = @{
    2  = "mardi"
    4  = "vendredi"
    7  = "mardi"
    9  = "vendredi"
    12 = "mardi"
    15 = "mardi"
    18 = "mardi"
    21 = "mardi"
    24 = "mardi"
    27 = "mardi"
    30 = "mardi"
    33 = "mardi"
}

foreach ($row in $dayRows.Keys) {
    $ws.Cells.Item($row, 2).Value = $dayRows[$row]
}
